$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8:G11").NumberFormat = "@"

# Row 8: assetId "1" (string), same rest pattern as row 4
$ws.Range("A8").Value = "1"
$ws.Range("B8").Value = "9860"
$ws.Range("C8").Value = "1234"
$ws.Range("D8").Value = "Welcome to St.Loius Airport "
$ws.Range("E8").Value = "123.5"
$ws.Range("F8").Value = "206"
$ws.Range("G8").Value = "St.Loius Airport shuttle 9"

# Row 9: assetId "2", same rest pattern as row 5
$ws.Range("A9").Value = "2"
$ws.Range("B9").Value = "9860"
$ws.Range("C9").Value = "5678"
$ws.Range("D9").Value = "Welcome to St.Loius Airport gateway 2"
$ws.Range("E9").Value = "109.5"
$ws.Range("F9").Value = "200"
$ws.Range("G9").Value = "St.Loius Airport shuttle 10"

# Row 10: assetId "123", same rest pattern as row 4
$ws.Range("A10").Value = "123"
$ws.Range("B10").Value = "9860"
$ws.Range("C10").Value = "1234"
$ws.Range("D10").Value = "Welcome to St.Loius Airport "
$ws.Range("E10").Value = "123.5"
$ws.Range("F10").Value = "206"
$ws.Range("G10").Value = "St.Loius Airport shuttle 9"

# Row 11: assetId "2", same rest pattern as row 5
$ws.Range("A11").Value = "2"
$ws.Range("B11").Value = "9860"
$ws.Range("C11").Value = "5678"
$ws.Range("D11").Value = "Welcome to St.Loius Airport gateway 2"
$ws.Range("E11").Value = "109.5"
$ws.Range("F11").Value = "200"
$ws.Range("G11").Value = "St.Loius Airport shuttle 10"

$ws.Range("A8:G11").Style = "Normal"

$ws.Range("F12").Select()
